$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.360.34"
$ws.Range("E2").Value = "  +0.42%  "

# Row 3
$ws.Range("D3").Value = "1.867.06"
$ws.Range("E3").Value = "  +0.24%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.86%  "

# Row 6
$ws.Range("E6").Value = "  +0.00%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4709"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.36%  "

# Row 8
$ws.Range("E8").Value = "  -1.10%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06564"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.40%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.10%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07878"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.76%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.95"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.25%  "

# Row 13
$ws.Range("D13").Value = "1.865.67"
$ws.Range("E13").Value = "  +0.13%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6931"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.80%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.105"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.20%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "

# Row 17
$ws.Range("D17").Value = "30.299.12"
$ws.Range("E17").Value = "  +0.25%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.04%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007648"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.60%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.03%  "

# Row 21
$ws.Range("D21").Value = "2.111.71"
$ws.Range("E21").Value = "  -0.12%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.229"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.71%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.178"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.407"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.33%  "

# Row 27
$ws.Range("E27").Value = "  +0.01%  "

# Row 28
$ws.Range("E28").Value = "  -0.68%  "

# Row 29
$ws.Range("E29").Value = "  -2.09%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09915"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.79%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.373"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.05%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.459"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.053"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04751"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.99%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.134"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.38%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7025"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.15%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.717"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.32%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01875"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.795"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.96%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.310"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.30%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.55%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.951"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.76%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4180"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.43%  "

# Row 44
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8426"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.38%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.07%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "969.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.53%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.123"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.39%  "

# Row 49
$ws.Range("E49").Value = "  -0.28%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.15%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05677"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.47%  "
